{"js": "// Neutralize racial-misclassification language in this (non-electoral) resume:\n//   \"affecting all Black and Asian-American voters\" -> \"affecting 50M voters\"\n// applied to the professional summary, the work-experience bullet (where the\n// replacement text is emphasized as bold/colored, matching the existing\n// \"50M\"-style metric runs in that bullet), and the project \"Impact:\" line.\n\nconst body = context.document.body;\n\n// 1) Professional summary paragraph \u2014 plain text swap (single run).\nconst summaryMatches = body.search(\n  \"affecting all Black and Asian-American voters, developed geospatial ML\",\n  { matchCase: true }\n);\nsummaryMatches.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < summaryMatches.items.length; i++) {\n  summaryMatches.items[i].insertText(\n    \"affecting 50M voters, developed geospatial ML\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n\n// 2) Work-experience bullet \u2014 \"50M\" becomes its own bold/colored run, just\n//    like the \"23%\"/\"64%\" metric runs later in the same bullet.\nconst bulletMatches = body.search(\n  \"all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from\",\n  { matchCase: true }\n);\nbulletMatches.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < bulletMatches.items.length; i++) {\n  bulletMatches.items[i].insertText(\n    \"50M voters, developed geospatial machine learning algorithms improving demographic classification accuracy from\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n\nconst boldMatches = body.search(\n  \"\u2022 Discovered systematic race coding errors affecting 50M\",\n  { matchCase: true }\n);\nboldMatches.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < boldMatches.items.length; i++) {\n  // Narrow the matched range down to just the \"50M\" token so only that\n  // portion gets split into its own bold/colored run.\n  const fiftyM = boldMatches.items[i].search(\"50M\", { matchCase: true });\n  fiftyM.load(\"items\");\n  await context.sync();\n  for (let j = 0; j < fiftyM.items.length; j++) {\n    fiftyM.items[j].font.bold = true;\n    fiftyM.items[j].font.color = \"#2C3E50\";\n  }\n}\nawait context.sync();\n\n// 3) Project \"Impact:\" line \u2014 plain text swap (single run), note the added\n//    \"nationwide\".\nconst impactMatches = body.search(\n  \"affecting all Black and Asian-American voters, improved electoral prediction accuracy\",\n  { matchCase: true }\n);\nimpactMatches.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < impactMatches.items.length; i++) {\n  impactMatches.items[i].insertText(\n    \"affecting 50M voters nationwide, improved electoral prediction accuracy\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n", "ps1": "# Neutralize racial-misclassification language in this (non-electoral) resume:\n#   \"affecting all Black and Asian-American voters\" -> \"affecting 50M voters\"\n# applied to the professional summary, the work-experience bullet (where the\n# replacement text is emphasized as bold/colored, matching the existing\n# \"50M\"-style metric runs in that bullet), and the project \"Impact:\" line.\n#\n# Each target paragraph is located by its distinctive wording (rather than a\n# fixed paragraph index) and every Find/Replace below is scoped to that\n# single paragraph's Range, so a \"50M\" introduced by one edit can never be\n# picked up by a later, unrelated search.\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n$count = $paras.Count\n\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $paras.Item($i)\n    $t = $p.Range.Text\n\n    if ($t -like \"*Product-focused data scientist*affecting all Black and Asian-American voters*\") {\n        # 1) Professional summary paragraph - plain text swap (single run).\n        $range = $p.Range\n        $find = $range.Find\n        $find.Text = \"affecting all Black and Asian-American voters, developed geospatial ML\"\n        $find.MatchCase = $true\n        $result = $find.Execute(\"affecting all Black and Asian-American voters, developed geospatial ML\", $true, $false, $false, $false, $false, $true, 1, $false, \"affecting 50M voters, developed geospatial ML\", 2)\n    }\n    elseif ($t -like \"*Discovered systematic race coding errors affecting all Black and Asian-American voters*\") {\n        # 2) Work-experience bullet - plain text swap first ...\n        $range = $p.Range\n        $find = $range.Find\n        $find.Text = \"all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from\"\n        $find.MatchCase = $true\n        $result = $find.Execute(\"all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from\", $true, $false, $false, $false, $false, $true, 1, $false, \"50M voters, developed geospatial machine learning algorithms improving demographic classification accuracy from\", 2)\n\n        # ... then split \"50M\" into its own bold/colored run, matching the\n        # other bold metric runs (\"23%\", \"64%\", etc.) already present in\n        # this bullet. Re-scope to the (now-changed) paragraph range first.\n        $range2 = $p.Range\n        $find2 = $range2.Find\n        $find2.Text = \"50M\"\n        $find2.MatchCase = $true\n        $result2 = $find2.Execute()\n        $range2.Font.Bold = 1\n        # w:color val=\"2C3E50\" -- COM Font.Color is a BGR-packed OLE_COLOR,\n        # so swap the RGB byte order before assigning.\n        $r = 0x2C\n        $g = 0x3E\n        $b = 0x50\n        $range2.Font.Color = ($b * 65536) + ($g * 256) + $r\n    }\n    elseif ($t -like \"*Impact: Corrected demographic data affecting all Black and Asian-American voters*\") {\n        # 3) Project \"Impact:\" line - plain text swap (single run); note the\n        #    added \"nationwide\".\n        $range = $p.Range\n        $find = $range.Find\n        $find.Text = \"affecting all Black and Asian-American voters, improved electoral prediction accuracy\"\n        $find.MatchCase = $true\n        $result = $find.Execute(\"affecting all Black and Asian-American voters, improved electoral prediction accuracy\", $true, $false, $false, $false, $false, $true, 1, $false, \"affecting 50M voters nationwide, improved electoral prediction accuracy\", 2)\n    }\n}\n"}
